$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 55505.05
$ws.Range("J64").Value = 3103
$ws.Range("L64").Value = 3103
$ws.Range("N64").Value = -3599

$ws.Range("H67").Value = 55505.05
$ws.Range("J67").Value = 3103
$ws.Range("L67").Value = 3103
$ws.Range("N67").Value = -4819

$ws.Range("H112").Value = 1033.2903
$ws.Range("J112").Value = 1051.0667
$ws.Range("L112").Value = 3153.2001
$ws.Range("N112").Value = -5369.2001

$ws.Range("H116").Value = 4393.6924
$ws.Range("I116").Value = 4593.1665
$ws.Range("K116").Value = 4593.1665
$ws.Range("M116").Value = -1151.1665

$ws.Range("H132").Value = 4722414
$ws.Range("I132").Value = 5005552.5
$ws.Range("J132").Value = 3433.3333
$ws.Range("K132").Value = 15016657.5
$ws.Range("L132").Value = 10299.9999
$ws.Range("M132").Value = -15014127.5
$ws.Range("N132").Value = -15359.9999

$ws.Range("H141").Value = 2692.889
$ws.Range("I141").Value = 2372.0715
$ws.Range("J141").Value = 3815.75
$ws.Range("K141").Value = 7116.2145
$ws.Range("L141").Value = 11447.25
$ws.Range("M141").Value = -1936.2145
$ws.Range("N141").Value = -21807.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1600.86
$ws.Range("I32").Value = 1600.86
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1600.86
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1313.86
$ws.Range("N32").ClearContents()

$ws.Range("H61").Value = 2100.4883
$ws.Range("I61").Value = 1104.95
$ws.Range("J61").Value = 2966.1738
$ws.Range("K61").Value = 1104.95
$ws.Range("L61").Value = 2966.1738
$ws.Range("M61").Value = -892.95
$ws.Range("N61").Value = -3390.1738

$ws.Range("H136").Value = 2100.4883
$ws.Range("I136").Value = 1104.95
$ws.Range("J136").Value = 2966.1738
$ws.Range("K136").Value = 3314.85
$ws.Range("L136").Value = 8898.5214
$ws.Range("M136").Value = -764.8500000000004
$ws.Range("N136").Value = -13998.5214

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 183931.9
$ws.Range("I105").Value = 168678.33
$ws.Range("J105").Value = 202236.2
$ws.Range("K105").Value = 168678.33
$ws.Range("L105").Value = 202236.2
$ws.Range("M105").Value = -166931.33
$ws.Range("N105").Value = -205730.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2072.923
$ws.Range("I58").Value = 1698.1177
$ws.Range("K58").Value = 1698.1177
$ws.Range("M58").Value = -1495.1177

$ws.Range("H94").Value = 1337.4375
$ws.Range("I94").Value = 1337.3334
$ws.Range("J94").Value = 1337.4615
$ws.Range("K94").Value = 1337.3334
$ws.Range("L94").Value = 1337.4615
$ws.Range("M94").Value = -886.3334
$ws.Range("N94").Value = -2239.4615

$ws.Range("H132").Value = 2305.6316
$ws.Range("I132").Value = 2472.8965
$ws.Range("J132").Value = 1766.6666
$ws.Range("K132").Value = 7418.689499999999
$ws.Range("L132").Value = 5299.9998
$ws.Range("M132").Value = -4888.689499999999
$ws.Range("N132").Value = -10359.9998

$ws.Range("H134").Value = 1228.2778
$ws.Range("I134").Value = 699.1539
$ws.Range("J134").Value = 2604
$ws.Range("K134").Value = 2097.4617
$ws.Range("L134").Value = 7812
$ws.Range("M134").Value = 437.5383000000002
$ws.Range("N134").Value = -12882

$ws.Range("H136").Value = 2072.923
$ws.Range("I136").Value = 1698.1177
$ws.Range("K136").Value = 5094.3531
$ws.Range("M136").Value = -2544.3531

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 2100
$ws.Range("J58").Value = 1400
$ws.Range("L58").Value = 4200
$ws.Range("N58").Value = -4456

$ws.Range("H131").Value = 852.51044
$ws.Range("J131").Value = 856.81915
$ws.Range("L131").Value = 2570.45745
$ws.Range("N131").Value = -12650.45745

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 1990
$ws.Range("I44").Value = 1990
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 1990
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -1394
$ws.Range("N44").ClearContents()

$ws.Range("H70").Value = 58971.324
$ws.Range("I70").Value = 75457.21000000001
$ws.Range("J70").Value = 7681.8887
$ws.Range("K70").Value = 75457.21000000001
$ws.Range("L70").Value = 7681.8887
$ws.Range("M70").Value = -75187.21000000001
$ws.Range("N70").Value = -8221.8887

$ws.Range("H73").Value = 58971.324
$ws.Range("I73").Value = 75457.21000000001
$ws.Range("J73").Value = 7681.8887
$ws.Range("K73").Value = 75457.21000000001
$ws.Range("L73").Value = 7681.8887
$ws.Range("M73").Value = -74521.21000000001
$ws.Range("N73").Value = -9553.8887

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1975.2142
$ws.Range("I22").Value = 2687.25
$ws.Range("J22").Value = 1690.4
$ws.Range("K22").Value = 2687.25
$ws.Range("L22").Value = 1690.4
$ws.Range("M22").Value = -2392.25
$ws.Range("N22").Value = -2280.4

$ws.Range("H27").Value = 1975.2142
$ws.Range("I27").Value = 2687.25
$ws.Range("J27").Value = 1690.4
$ws.Range("K27").Value = 2687.25
$ws.Range("L27").Value = 1690.4
$ws.Range("M27").Value = -2580.25
$ws.Range("N27").Value = -1904.4

$ws.Range("H40").Value = 55124.26
$ws.Range("I40").Value = 201738.6
$ws.Range("J40").Value = 2762
$ws.Range("K40").Value = 201738.6
$ws.Range("L40").Value = 2762
$ws.Range("M40").Value = -201602.6
$ws.Range("N40").Value = -3034

$ws.Range("H68").Value = 5007.143
$ws.Range("J68").Value = 5558.3335
$ws.Range("L68").Value = 5558.3335
$ws.Range("N68").Value = -7056.3335

$ws.Range("H71").Value = 5007.143
$ws.Range("J71").Value = 5558.3335
$ws.Range("L71").Value = 27791.6675
$ws.Range("N71").Value = -35279.6675

$ws.Range("H122").Value = 2788.0557
$ws.Range("I122").Value = 2605.8
$ws.Range("K122").Value = 7817.400000000001
$ws.Range("M122").Value = -5367.400000000001

$ws.Range("H127").Value = 38000
$ws.Range("J127").Value = 38000
$ws.Range("L127").Value = 38000
$ws.Range("N127").Value = -47920

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9617559
$ws.Range("I62").Value = 15386474
$ws.Range("K62").Value = 15386474
$ws.Range("M62").Value = -15385850

$ws.Range("H65").Value = 9617559
$ws.Range("I65").Value = 15386474
$ws.Range("K65").Value = 76932370
$ws.Range("M65").Value = -76929250

$ws.Range("H129").Value = 32501.8
$ws.Range("J129").Value = 32501.8
$ws.Range("L129").Value = 32501.8
$ws.Range("N129").Value = -42501.8

$ws.Range("H132").Value = 3330.25
$ws.Range("I132").Value = 3862.476
$ws.Range("J132").Value = 1733.5714
$ws.Range("K132").Value = 11587.428
$ws.Range("L132").Value = 5200.7142
$ws.Range("M132").Value = -9057.428
$ws.Range("N132").Value = -10260.7142

$ws.Range("H136").Value = 1707.6111
$ws.Range("I136").Value = 705.087
$ws.Range("K136").Value = 2115.261
$ws.Range("M136").Value = 434.739
